$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
  "https://www.uyyaala.com/products/nestle-pre-nan-for-premature-babies-low-birth-weight-infant-formula-less-than-2-5kg-400g",
  "https://www.uyyaala.com/products/nestle-nan-optipro-starter-infant-formula-stage-1-400g-0-to-6months",
  "https://www.uyyaala.com/products/products-nestle-nan-optipro-starter-follow-up-formula-stage-3-400g-1-to-3-years",
  "https://www.uyyaala.com/products/nestle-lactogen-infant-formula-stage-2-after-6-months-400g",
  "https://www.uyyaala.com/products/nestle-nan-pro-infant-formula-stage-1-400g-upto-6-months",
  "https://www.uyyaala.com/products/nestle-lactogen-infant-formula-stage-1-up-to-6-months-400g-tin-pack",
  "https://www.uyyaala.com/products/nestle-nan-pro-infant-formula-stage-2-after-6-months-400g",
  "https://www.uyyaala.com/products/nestle-nido-one-plus-growing-up-formula-1-3-years-400g",
  "https://www.uyyaala.com/products/nestle-lactogen-infant-formula-stage-1-upto-6-months-400g",
  "https://www.uyyaala.com/products/nestle-nan-pro-infant-formula-stage-3-after-12-months-400g",
  "https://www.uyyaala.com/products/nestle-lactogen-infant-formula-stage-3-after-12months-400g",
  "https://www.uyyaala.com/products/nestle-nan-optipro-starter-follow-up-formula-stage-4-400g-3-to-6-years",
  "https://www.uyyaala.com/products/nestle-lactogen-infant-formula-stage-4-18-to-24months-400g",
  "https://www.uyyaala.com/products/nestle-nan-pro-infant-baby-formula-stage-4-18-to-24-months-400g",
  "https://www.uyyaala.com/products/nestle-excellapro-infant-formula-stage-1-up-to-6-months-400g",
  "https://www.uyyaala.com/products/nestle-nido-little-kids-infant-formula-3-5-years-1800g"
)

# Remove the now-unused trailing rows entirely (rows 18-27) so the sheet
# dimension shrinks back down to A1:B17.
$ws.Rows("18:27").Delete() | Out-Null

# Write the new url list into A2:B17 (index in col A, url text in col B).
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $urls[$i]
}

$wb.Save()
